$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.086.74'
$ws.Range('E2').Value = '  +3.56%  '
$ws.Range('D3').Value = '1.916.85'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4985'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.92%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3001'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06900'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.11%  '
$ws.Range('D10').Value = '1.915.16'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07309'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '89.67'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6824'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.093'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.08%  '
$ws.Range('D16').Value = '31.030.05'
$ws.Range('E16').Value = '  +3.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008050'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.15%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '2.161.15'
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.880'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '175.79'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +29.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.088'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +9.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.336'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.949'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.412'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.356'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08960'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.062'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05257'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7500'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.85%  '
$ws.Range('E35').Value = '  +2.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.668'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01912'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +16.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.742'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.197'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9403'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.944'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4364'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '105.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.34%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.816'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1337'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.23%  '
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.600'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.68%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3892'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.392'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.48%  '
